# Weekly update: insert 4 new "Fruta" (Pera) price rows for Vega Monumental
# Concepción, dated 2021-09-09 (serial 44448), ahead of the existing data
# which shifts down from rows 225-234 to rows 229-238.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 225:234 down by four rows.
$ws.Rows("225:228").Insert()

$newRows = @(
    @{ Row=225; Variedad="Packham's Triumph"; Calidad="Primera"; Volumen=100; PMin=9000;  PMax=10000; PProm=9500; PKg=594 },
    @{ Row=226; Variedad="Packham's Triumph"; Calidad="Segunda"; Volumen=50;  PMin=8000;  PMax=8000;  PProm=8000; PKg=500 },
    @{ Row=227; Variedad="Winter Nelis";       Calidad="Primera"; Volumen=100; PMin=9000;  PMax=10000; PProm=9500; PKg=594 },
    @{ Row=228; Variedad="Winter Nelis";       Calidad="Segunda"; Volumen=50;  PMin=8000;  PMax=8000;  PProm=8000; PKg=500 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = 44448
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100104
    $ws.Cells.Item($row, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value = 100104005
    $ws.Cells.Item($row, 10).Value = "Pera"
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = "$/caja 16 kilos empedrada"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = 16
}
